$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper functions for working with the Hyperlinks collection: the
# collection mutates-while-iterating badly, so always re-scan from scratch
# and stop after the first match.
# ---------------------------------------------------------------------------
function Remove-HyperlinkAt($sheet, $targetAddr) {
    foreach ($hl in $sheet.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq $targetAddr) {
            $hl.Delete()
            break
        }
    }
}

function Set-HyperlinkAddress($sheet, $targetAddr, $newAddr) {
    foreach ($hl in $sheet.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq $targetAddr) {
            $hl.Address = $newAddr
            break
        }
    }
}

# ---------------------------------------------------------------------------
# Row 1 - header: "Banyak" -> "Qy"
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Qy"

# ---------------------------------------------------------------------------
# Row 2 - Module Bluetooth HC-05: price bump 43000 -> 43500
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 43500
$ws.Range("D2").Value = 43500

# ---------------------------------------------------------------------------
# Row 3 - Arduino Nano/Uno -> Arduino Uno, qty 1 -> 2, add "Toko Offline" note
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "Toko Offline"
$ws.Range("A3").Value = "Arduino Uno"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 150000

# ---------------------------------------------------------------------------
# Row 4 - Dinamo -> Kit Chassis Car 4 WD, price 17000 -> 170000, qty 4 -> 1
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Kit Chassis Car 4 WD"
$ws.Range("B4").Value = 170000
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 170000
$ws.Range("E4").Value = "iSee Surabaya"

# ---------------------------------------------------------------------------
# Row 5 - L293 Motor Driver -> L293D Shield, price 27000 -> 35000
# reference note changes to "iSee Surabaya" and the mailto hyperlink is removed
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "L293D Shield"
$ws.Range("B5").Value = 35000
$ws.Range("D5").Value = 35000
$ws.Range("E5").Value = "iSee Surabaya"
Remove-HyperlinkAt $ws '$E$5'

# ---------------------------------------------------------------------------
# Row 6 - Sensor Ultrasonic HC-SR04 -> Sensor Ultrasonic , price 14500 -> 10000
# qty 2 -> 4
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Sensor Ultrasonic "
$ws.Range("B6").Value = 10000
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 40000

# ---------------------------------------------------------------------------
# Row 7 - Chasiss Car 4WD -> Battery Holder  3x 18650, price 105000 -> 10500
# reference text/hyperlink switches from Shopee/@nasrula to Shopee/@Mulia-Electric
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Battery Holder  3x 18650"
$ws.Range("B7").Value = 10500
$ws.Range("D7").Value = 10500
$ws.Range("E7").Value = "Shopee/@Mulia-Electric"
Set-HyperlinkAddress $ws '$E$7' "mailto:Shopee/@Mulia-Electric"

# ---------------------------------------------------------------------------
# Row 8 - "Holder Baterai 18650" item removed entirely; clear the row and
# drop its hyperlink, keeping the existing formatting in place.
# ---------------------------------------------------------------------------
$ws.Range("A8:E8").ClearContents()
Remove-HyperlinkAt $ws '$E$8'

# ---------------------------------------------------------------------------
# Row 9 - new grand-total formula (moved up from row 10), currency style
# ---------------------------------------------------------------------------
$ws.Range("D9").Style = "Currency [0]"
$ws.Range("D9").NumberFormat = '_-"Rp"* #,##0_-;\-"Rp"* #,##0_-;_-"Rp"* "-"_-;_-@_-'
$ws.Range("D9").Formula = "=SUM(D2:D7)"

# ---------------------------------------------------------------------------
# Row 10 - old grand-total formula & label cleared (kept blank, same styles)
# ---------------------------------------------------------------------------
$ws.Range("D10:E10").ClearContents()

# ---------------------------------------------------------------------------
# Selection, to match the author's last-active cell
# ---------------------------------------------------------------------------
$ws.Range("E12").Select() | Out-Null
